$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.72856851286642
$ws.Range("C2").Value = 9.916638121763459
$ws.Range("D2").Value = 3.505671535895402
$ws.Range("E2").Value = 16.59418520955057
$ws.Range("F2").Value = 20.07054344953894
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 16.4614748518014
$ws.Range("N2").Value = 15.59264994587608
$ws.Range("O2").Value = 17.64544065876387
$ws.Range("B3").Value = 12.05561713470845
$ws.Range("C3").Value = 9.359310192563211
$ws.Range("D3").Value = 3.468360177280706
$ws.Range("E3").Value = 15.64598135373817
$ws.Range("F3").Value = 19.95179175234761
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 16.55071082291001
$ws.Range("N3").Value = 15.63248803818554
$ws.Range("O3").Value = 17.63198918109224
$ws.Range("B4").Value = 11.62337329251245
$ws.Range("C4").Value = 8.997907736710186
$ws.Range("D4").Value = 3.444993659384588
$ws.Range("E4").Value = 15.03844741167595
$ws.Range("F4").Value = 19.88636340653288
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 16.61033915098415
$ws.Range("N4").Value = 15.65871563750926
$ws.Range("O4").Value = 17.62987676238334
$ws.Range("B5").Value = 11.44261975362477
$ws.Range("C5").Value = 8.845862054095594
$ws.Range("D5").Value = 3.435361996828866
$ws.Range("E5").Value = 14.78476936934786
$ws.Range("F5").Value = 19.86160635224642
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 16.63584863546339
$ws.Range("N5").Value = 15.66984858985222
$ws.Range("O5").Value = 17.63055894528079
$ws.Range("B6").Value = 11.41233319752832
$ws.Range("C6").Value = 8.820328754983253
$ws.Range("D6").Value = 3.433756204080175
$ws.Range("E6").Value = 14.74228657337714
$ws.Range("F6").Value = 19.85761115384798
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 16.64015740362417
$ws.Range("N6").Value = 15.67172410556527
$ws.Range("O6").Value = 17.63076530716177
$ws.Range("B7").Value = 11.62095399396892
$ws.Range("C7").Value = 8.99587642857845
$ws.Range("D7").Value = 3.444864200138147
$ws.Range("E7").Value = 15.03505053970443
$ws.Range("F7").Value = 19.88602178111449
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 16.61067828750961
$ws.Range("N7").Value = 15.65886397771386
$ws.Range("O7").Value = 17.62987971969188
$ws.Range("B8").Value = 12.50058614978716
$ws.Range("C8").Value = 9.728486622238446
$ws.Range("D8").Value = 3.49290464630142
$ws.Range("E8").Value = 16.27264693760583
$ws.Range("F8").Value = 20.02805929256505
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 16.4912359439123
$ws.Range("N8").Value = 15.60601982957008
$ws.Range("O8").Value = 17.63952567459267
$ws.Range("B9").Value = 14.06805789349014
$ws.Range("C9").Value = 11.01091885283069
$ws.Range("D9").Value = 3.583255378082351
$ws.Range("E9").Value = 18.60726830698778
$ws.Range("F9").Value = 20.36478769375605
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 16.2956533549361
$ws.Range("N9").Value = 15.51638103031688
$ws.Range("O9").Value = 17.70726183208841
$ws.Range("B10").Value = 15.11697594969355
$ws.Range("C10").Value = 11.85733233898423
$ws.Range("D10").Value = 3.646988075677346
$ws.Range("E10").Value = 20.26377836507915
$ws.Range("F10").Value = 20.64580116530065
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 16.17587901396221
$ws.Range("N10").Value = 15.45901004345274
$ws.Range("O10").Value = 17.78675238888261
$ws.Range("B11").Value = 15.57085862463523
$ws.Range("C11").Value = 12.22137222447142
$ws.Range("D11").Value = 3.675346414604514
$ws.Range("E11").Value = 20.97493147040734
$ws.Range("F11").Value = 20.78047949586021
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 16.12666288716441
$ws.Range("N11").Value = 15.43474496541171
$ws.Range("O11").Value = 17.82932700121259
$ws.Range("B12").Value = 15.73932325125847
$ws.Range("C12").Value = 12.35619419434914
$ws.Range("D12").Value = 3.685988920103151
$ws.Range("E12").Value = 21.23816649240475
$ws.Range("F12").Value = 20.83242000607743
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 16.10879035404277
$ws.Range("N12").Value = 15.42581942786339
$ws.Range("O12").Value = 17.84636498338837
$ws.Range("B13").Value = 15.70319395860103
$ws.Range("C13").Value = 12.32729287232737
$ws.Range("D13").Value = 3.68370122193023
$ws.Range("E13").Value = 21.18174324382001
$ws.Range("F13").Value = 20.8211925626998
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 16.11260541282912
$ws.Range("N13").Value = 15.42773000746672
$ws.Range("O13").Value = 17.8426549336916
$ws.Range("B14").Value = 15.58478698954652
$ws.Range("C14").Value = 12.23252493100999
$ws.Range("D14").Value = 3.676223935706929
$ws.Range("E14").Value = 20.9967092957434
$ws.Range("F14").Value = 20.78473407549773
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 16.12517714176581
$ws.Range("N14").Value = 15.43400538553325
$ws.Range("O14").Value = 17.8307104128204
$ws.Range("B15").Value = 15.51181343239581
$ws.Range("C15").Value = 12.1740816772202
$ws.Range("D15").Value = 3.671631210228449
$ws.Range("E15").Value = 20.88258205652712
$ws.Range("F15").Value = 20.76252338583522
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 16.13297745572845
$ws.Range("N15").Value = 15.43788348736108
$ws.Range("O15").Value = 17.82351310778245
$ws.Range("B16").Value = 15.08683951609671
$ws.Range("C16").Value = 11.8331177619222
$ws.Range("D16").Value = 3.645121601617081
$ws.Range("E16").Value = 20.2164525261526
$ws.Range("F16").Value = 20.63713368015328
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 16.17920202282549
$ws.Range("N16").Value = 15.46063267016713
$ws.Range("O16").Value = 17.78409859178338
$ws.Range("B17").Value = 14.82012054354453
$ws.Range("C17").Value = 11.618558893633
$ws.Range("D17").Value = 3.628692816015363
$ws.Range("E17").Value = 19.79696787864455
$ws.Range("F17").Value = 20.5619339222156
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 16.20891382595779
$ws.Range("N17").Value = 15.47505773878762
$ws.Range("O17").Value = 17.76155778596574
$ws.Range("B18").Value = 14.66452292684296
$ws.Range("C18").Value = 11.49317430154187
$ws.Range("D18").Value = 3.619184029171352
$ws.Range("E18").Value = 19.55169833973925
$ws.Range("F18").Value = 20.51932794794511
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 16.22649888492742
$ws.Range("N18").Value = 15.48352724155493
$ws.Range("O18").Value = 17.7491971119968
$ws.Range("B19").Value = 14.6114666413611
$ws.Range("C19").Value = 11.45038231440792
$ws.Range("D19").Value = 3.61595447362586
$ws.Range("E19").Value = 19.46796743330377
$ws.Range("F19").Value = 20.50501469291886
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 16.23253780072185
$ws.Range("N19").Value = 15.48642452690525
$ws.Range("O19").Value = 17.74511594113448
$ws.Range("B20").Value = 14.84874018020895
$ws.Range("C20").Value = 11.64160366668036
$ws.Range("D20").Value = 3.630447871135758
$ws.Range("E20").Value = 19.84203575524988
$ws.Range("F20").Value = 20.56987241612021
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 16.20569961049588
$ws.Range("N20").Value = 15.47350430739908
$ws.Range("O20").Value = 17.76389480188505
$ws.Range("B21").Value = 15.61965900133971
$ws.Range("C21").Value = 12.26044293820684
$ws.Range("D21").Value = 3.678422845421021
$ws.Range("E21").Value = 21.05122253528351
$ws.Range("F21").Value = 20.79541764235748
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 16.12146371682025
$ws.Range("N21").Value = 15.43215501774795
$ws.Range("O21").Value = 17.83419401010477
$ws.Range("B22").Value = 16.10359692593644
$ws.Range("C22").Value = 12.64721444363892
$ws.Range("D22").Value = 3.709214434573908
$ws.Range("E22").Value = 21.80617053841991
$ws.Range("F22").Value = 20.94828624192638
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 16.07087018670322
$ws.Range("N22").Value = 15.40666431795305
$ws.Range("O22").Value = 17.8854730544125
$ws.Range("B23").Value = 15.8471490004155
$ws.Range("C23").Value = 12.44240752934946
$ws.Range("D23").Value = 3.692833496619246
$ws.Range("E23").Value = 21.40646166072641
$ws.Range("F23").Value = 20.86621254688848
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 16.09746262925452
$ws.Range("N23").Value = 15.42012903379739
$ws.Range("O23").Value = 17.85761890148904
$ws.Range("B24").Value = 14.83580826854918
$ws.Range("C24").Value = 11.63119145264961
$ws.Range("D24").Value = 3.629654608588755
$ws.Range("E24").Value = 19.82167334911489
$ws.Range("F24").Value = 20.56628146544555
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 16.2071511898234
$ws.Range("N24").Value = 15.47420606437038
$ws.Range("O24").Value = 17.76283637263846
$ws.Range("B25").Value = 13.66165326652721
$ws.Range("C25").Value = 10.68071316630868
$ws.Range("D25").Value = 3.559257126963713
$ws.Range("E25").Value = 17.95941609507481
$ws.Range("F25").Value = 20.26764239801417
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 16.3443893052115
$ws.Range("N25").Value = 15.53913746813186
$ws.Range("O25").Value = 17.68370822244885
